$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 7014
$ws.Range("E7").Value = 290682013

$ws.Range("C37").Value = 23050
$ws.Range("E37").Value = 130197201

$ws.Range("C51").Value = 6357
$ws.Range("E51").Value = 12201689

$ws.Range("C92").Value = 409282
$ws.Range("E92").Value = 1597152985

$ws.Range("C93").Value = 209656
$ws.Range("E93").Value = 1310025277

$ws.Range("C94").Value = 94231
$ws.Range("E94").Value = 919017349

$ws.Range("C95").Value = 50802
$ws.Range("E95").Value = 934242292

$ws.Range("C116").Value = 4566
$ws.Range("E116").Value = 20667158

$ws.Range("C121").Value = 14
$ws.Range("E121").Value = 1153896
